$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header fields
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Statement header
$ws.Range("D5").Value = "KONTOSTAND AM 23.10.2024"

# Row 6
$ws.Range("B6").Value = "26.10."
$ws.Range("C6").Value = "27.10."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 23544901"
$ws.Range("E6").Value = "38,25-"

# Row 7
$ws.Range("B7").Value = "30.10."
$ws.Range("C7").Value = "31.10."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "25,16-"

# Row 8
$ws.Range("B8").Value = "01.11."
$ws.Range("C8").Value = "02.11."
$ws.Range("D8").Value = "MCDONALDS Borna"
$ws.Range("E8").Value = "30,51-"

# Row 9 (previously empty) - copy E8's formatting to E9 (s changes 13 -> 17)
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("B9").Value = "03.11."
$ws.Range("C9").Value = "04.11."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-2656831"
$ws.Range("E9").Value = "54,67-"

# Row 10 (previously empty) - copy E8's formatting to E10 (s changes 12 -> 17)
$ws.Range("E8").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("B10").Value = "05.11."
$ws.Range("C10").Value = "06.11."
$ws.Range("D10").Value = "AMAZON.DE MKTPLC EU MZLMPB"
$ws.Range("E10").Value = "125,72-"

# Final balance / next billing date
$ws.Range("D12").Value = "KONTOSTAND AM 10.11.2024"
$ws.Range("E12").Value = "274,31-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 17.11.2024"
